$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.172.11"
$ws.Range("E2").Value = "  -4.92%  "
$ws.Range("D3").Value = "3.001.80"
$ws.Range("E3").Value = "  -5.14%  "
$ws.Range("E4").Value = "  -0.11%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "570.20"
$c.ClearFormats()
$ws.Range("E5").Value = "  -4.77%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "125.33"
$c.ClearFormats()
$ws.Range("E6").Value = "  -7.31%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.003.86"
$ws.Range("E8").Value = "  -4.94%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.504"
$c.ClearFormats()
$ws.Range("E9").Value = "  -2.16%  "
$ws.Range("E10").Value = "  -6.92%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.07"
$c.ClearFormats()
$ws.Range("E11").Value = "  -4.48%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.442"
$c.ClearFormats()
$ws.Range("E12").Value = "  -2.47%  "
$ws.Range("E13").Value = "  -7.23%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "32.55"
$c.ClearFormats()
$ws.Range("E14").Value = "  -5.57%  "
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").Value = "3.494.73"
$ws.Range("E16").Value = "  -5.23%  "
$ws.Range("D17").Value = "2.994.61"
$ws.Range("E17").Value = "  -5.49%  "
$ws.Range("D18").Value = "60.108.79"
$ws.Range("E18").Value = "  -5.10%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.52"
$c.ClearFormats()
$ws.Range("E19").Value = "  -0.41%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "430.21"
$c.ClearFormats()
$ws.Range("E20").Value = "  -6.48%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.16"
$c.ClearFormats()
$ws.Range("E21").Value = "  -5.45%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.673"
$c.ClearFormats()
$ws.Range("E22").Value = "  -2.93%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.08"
$c.ClearFormats()
$ws.Range("E23").Value = "  -7.19%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "12.88"
$c.ClearFormats()
$ws.Range("E24").Value = "  -2.15%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "79.60"
$c.ClearFormats()
$ws.Range("E25").Value = "  -4.09%  "
$ws.Range("E26").Value = "  +0.18%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("E28").Value = "  -6.22%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.97"
$c.ClearFormats()
$ws.Range("E29").Value = "  -4.46%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.21"
$c.ClearFormats()
$ws.Range("E30").Value = "  -6.33%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "6.13"
$c.ClearFormats()
$ws.Range("E31").Value = "  -9.82%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "25.31"
$c.ClearFormats()
$ws.Range("E32").Value = "  -6.85%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0950"
$c.ClearFormats()
$ws.Range("E33").Value = "  -5.43%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.61"
$c.ClearFormats()
$ws.Range("E34").Value = "  -4.29%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.935"
$c.ClearFormats()
$ws.Range("E35").Value = "  -8.22%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "50.37"
$c.ClearFormats()
$ws.Range("E36").Value = "  -1.73%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.03"
$c.ClearFormats()
$ws.Range("E37").Value = "  -15.77%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "8.50"
$c.ClearFormats()
$ws.Range("E38").Value = "  +4.81%  "
$ws.Range("D39").Value = "0.0₃0662"
$ws.Range("E39").Value = "  -9.71%  "
$ws.Range("E40").Value = "  -8.32%  "
$ws.Range("E41").Value = "  -3.99%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "372.15"
$c.ClearFormats()
$ws.Range("E42").Value = "  -4.52%  "
$ws.Range("D43").Value = "2.673.63"
$ws.Range("E43").Value = "  -4.01%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.44"
$c.ClearFormats()
$ws.Range("E44").Value = "  -6.71%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "121.36"
$c.ClearFormats()
$ws.Range("E46").Value = "  -4.74%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.234"
$c.ClearFormats()
$ws.Range("E47").Value = "  -6.39%  "
$ws.Range("E48").Value = "  -5.50%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.108"
$c.ClearFormats()
$ws.Range("E49").Value = "  -3.14%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "23.31"
$c.ClearFormats()
$ws.Range("E50").Value = "  -6.45%  "
$ws.Range("E51").Value = "  -2.45%  "
